# Update CPI workbook (InputData/cpi.xlsx) to the newest EPS-US vintage:
#  - "About" sheet: point the source citation at the 2022-12 BLS PDF and
#    relabel "Page 4" -> "Pages 4 and 5"; drop the old hyperlink.
#  - "Data" sheet: shift the 1968-2020 data down (index bookkeeping is
#    handled automatically by shared strings), then append the 2021 and
#    2022 rows of CPI-U data pulled from BLS, with the "multiple" column
#    formula carried down.

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# About sheet: refresh the citation, drop the stale hyperlink, move the
# active selection to reflect the editor's final cursor position.
# ---------------------------------------------------------------------
$wsAbout.Range("B6").Value = "https://www.bls.gov/cpi/tables/supplemental-files/historical-cpi-u-202212.pdf"
$wsAbout.Range("B7").Value = "Pages 4 and 5"

# The old hyperlink (B6 -> old BLS pdf) is removed; the cell keeps the
# Hyperlink-style formatting but is plain text going forward.
$wsAbout.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# Data sheet: append 2021 and 2022 rows under the existing historical
# CPI-U table (row 58 is the last existing row, for 2020).
# ---------------------------------------------------------------------
$wsData.Range("A59").Value = "2021............................................................................. ."
$wsData.Cells.Item(59, 2).Value = 266.236
$wsData.Cells.Item(59, 3).Value = 275.703
$wsData.Cells.Item(59, 4).Value = 270.97
$wsData.Cells.Item(59, 5).Value = 7
$wsData.Cells.Item(59, 6).Value = 4.7
$wsData.Range("G59").Formula = '=$D$50/D59'

$wsData.Range("A60").Value = "2022............................................................................. ."
$wsData.Cells.Item(60, 2).Value = 288.347
$wsData.Cells.Item(60, 3).Value = 296.963
$wsData.Cells.Item(60, 4).Value = 292.655
$wsData.Cells.Item(60, 5).Value = 6.5
$wsData.Cells.Item(60, 6).Value = 8
$wsData.Range("G60").Formula = '=$D$50/D60'

# Blank styled cell below the new rows (matches the sheet's habit of
# carrying the 0.000 number format one row past the data).
$wsData.Range("G62").NumberFormat = "0.000"

# ---------------------------------------------------------------------
# Update sheet selections/activation to match where the editor left off:
# the About tab is now the active tab, selected cell B8; the Data tab's
# selection moved to the new last row, scrolled down a few rows further.
# ---------------------------------------------------------------------
$wsData.Range("B63").Select()
$wsAbout.Activate()
$wsAbout.Range("B8").Select()
